# ---------------------------------------------------------------------------
# Applies the two changes captured in the commit:
#   1. Bump the "datetimeFigureOut" date field text shown on the slide
#      master and every slide layout from 21.09.2023 -> 26.09.2023.
#   2. Re-color every shape that used the old accent red (C00000) to the
#      new accent green (769E3C), including giving the one picture in the
#      "struktur" diagram (slide 3) a matching solid fill it didn't have
#      before.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- helper: VBA/COM RGB long values are packed as 0x00BBGGRR, i.e.
#     R + G*256 + B*65536 -- build the new accent-green value from its hex.
function Get-BgrLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$oldDate = "21.09.2023"
$newDate = "26.09.2023"

$oldFillRgb = Get-BgrLong "C00000"
$newFillRgb = Get-BgrLong "769E3C"

# --- 1. Update the date placeholder text wherever it still shows the old
#        cached value: the slide master and each of its custom layouts.
function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# --- 2. Re-color shapes/pictures across every slide.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)

        if ($shp.Type -eq 13) {
            # msoPicture: the "struktur" diagrams have exactly one picture
            # per slide; only the one on slide 3 gains a new solid fill.
            if ($shp.Name -eq "Grafik 17") {
                $shp.Fill.Visible = $true
                $shp.Fill.ForeColor.RGB = $newFillRgb
            }
        } else {
            $curRgb = $null
            try { $curRgb = $shp.Fill.ForeColor.RGB } catch {}
            if ($curRgb -eq $oldFillRgb) {
                $shp.Fill.ForeColor.RGB = $newFillRgb
            }
        }
    }
}
